$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the source date text in O2 (used by P2's formula)
$ws.Range("O2").Value = "07/09/2024"

# S2 becomes a formula that mirrors P2 instead of a static number
$ws.Range("S2").Formula = "=P2"

# Update the active selection: moved from O1:P2 to S3
$ws.Range("S3").Select()
